$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '264.40'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.80'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.182'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06118'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.534'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.730'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.368'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8155'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1591'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08192'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03372'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03159'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09258'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.915'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001701'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04843'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0006208'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006244'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001104'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.003336'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001503'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.700'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.265'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002685'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04646'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1125'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003136'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003457'
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006160'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.7510'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1651'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002103'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01242'
